$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 22, shifting existing rows 22-25 down to 23-26.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly record,
# mirroring the structure of the surrounding rows.
$ws.Cells.Item(22, 1).Value = 1
$ws.Cells.Item(22, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(22, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(22, 4).Value = 45147
$ws.Cells.Item(22, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(22, 5).Value = 15
$ws.Cells.Item(22, 6).Value = "Fruta"
$ws.Cells.Item(22, 7).Value = 100107
$ws.Cells.Item(22, 8).Value = "Otros"
$ws.Cells.Item(22, 9).Value = 100107002
$ws.Cells.Item(22, 10).Value = "Chirimoya"
$ws.Cells.Item(22, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(22, 12).Value = "Primera"
$ws.Cells.Item(22, 13).Value = 270
$ws.Cells.Item(22, 14).Value = 17000
$ws.Cells.Item(22, 15).Value = 18000
$ws.Cells.Item(22, 16).Value = 17500
$ws.Cells.Item(22, 17).Value = "$/caja 10 kilos"
$ws.Cells.Item(22, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(22, 19).Value = 1750
$ws.Cells.Item(22, 20).Value = 10
